$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsrapport_Maikoff")
$ws.Range("D11").Value = 9
$sh = $ws.Shapes.Item(1)
$chart = $sh.Chart
$sc = $chart.SeriesCollection
$ser = $sc.Item(1)
try {
  $ser.Formula = "=SERIES(Arbeitsrapport_Maikoff!`$F`$5,Arbeitsrapport_Maikoff!`$G`$4:`$AG`$4,Arbeitsrapport_Maikoff!`$G`$5:`$AG`$5,1)"
  Write-Output "set formula ok"
} catch {
  Write-Output ("ERROR: " + $_.Exception.Message)
}
